$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 10-18 (player, position, team)
$data = @(
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Davion Mitchell", "PG,SG", "Toronto Raptors"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Jordan Poole", "PG,SG", "Washington Wizards")
)

$startRow = 10
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Remove the now-obsolete last row (former row 19), shifting cells up
$ws.Rows.Item(19).Delete()
